$wb = $excel.ActiveWorkbook

# --- Insert the new "Sheet1" worksheet before "关卡设计" ------------------
$target = $wb.Worksheets.Item("关卡设计")
$new = $wb.Worksheets.Add($target, [System.Reflection.Missing]::Value)

# Fill content in the same order new unique strings were introduced so the
# shared-string table ends up in the same order as the authored workbook.
$new.Cells.Item(1, 1).Value = "推箱子"
$new.Cells.Item(2, 1).Value = "旋轉"
$new.Cells.Item(3, 1).Value = "傳送門"
$new.Cells.Item(6, 1).Value = "冰(温度)"
$new.Cells.Item(7, 1).Value = "聊天記錄"
$new.Cells.Item(8, 1).Value = "分數GPA"
$new.Cells.Item(4, 1).Value = "移動(單個/多個)"
$new.Cells.Item(5, 1).Value = "日月改變(温度)(明暗障礙物)"

$new.Columns.Item(1).ColumnWidth = 71.125
$new.Range("A1:A8").RowHeight = 22.2
$new.Range("A1:A8").Font.Size = 16

$new.Range("G10").Select()

# --- View-state tweaks on the other sheets --------------------------------
$sheet2 = $wb.Worksheets.Item("优化和添加功能")
$sheet2.Range("D19").Select()

$sheet4 = $wb.Worksheets.Item("遊戲機制")
$sheet4.Application.ActiveWindow.ScrollRow = 26
$sheet4.Range("B39").Select()

$sheet6 = $wb.Worksheets.Item("美术需求")
$sheet6.Range("D16").Select()
$sheet6.Application.ActiveWindow.ScrollRow = 10

$sheet7 = $wb.Worksheets.Item("系统策划")
$sheet7.Range("D15").Select()

$sheet8 = $wb.Worksheets.Item("关卡设计")
$sheet8.Range("B20").Select()

# Make the new sheet the active tab (matches activeTab=7 / tabSelected on Sheet1)
$new.Activate()
$new.Range("G10").Select()
